$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-reg_center_user_h")

# Add new row 33 with data following the same pattern as prior rows
$ws.Range("A33").Value = 10002
$ws.Range("B33").Value = 110032
$ws.Range("C33").Value = "eng"
$ws.Range("D33").Value = $true
$ws.Range("E33").Value = "superadmin"
$ws.Range("F33").Value = "now()"
$ws.Range("G33").Value = "now()"

# Update the selected cell / view to match the saved state
$ws.Range("D26").Select()
